$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "descripcion"
$ws.Range("D1").Value = "informacion"
